# Auto-generated Excel COM-interop script to apply leve price/profit refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2826.875
$ws.Range("J2").Value = 6988.6665
$ws.Range("L2").Value = 6988.6665
$ws.Range("N2").Value = -7214.6665
$ws.Range("H9").Value = 508.0909
$ws.Range("I9").Value = 593.2222
$ws.Range("K9").Value = 593.2222
$ws.Range("M9").Value = -424.2222
$ws.Range("H12").Value = 290.2
$ws.Range("J12").Value = 289
$ws.Range("L12").Value = 289
$ws.Range("N12").Value = -629
$ws.Range("H16").Value = 7849.5
$ws.Range("J16").Value = 7849.5
$ws.Range("L16").Value = 7849.5
$ws.Range("N16").Value = -8309.5
$ws.Range("H18").Value = 973.2941
$ws.Range("I18").Value = 897.44446
$ws.Range("K18").Value = 897.44446
$ws.Range("M18").Value = -613.44446
$ws.Range("H43").Value = 1435.6
$ws.Range("I43").Value = 1426.6666
$ws.Range("K43").Value = 1426.6666
$ws.Range("M43").Value = -1357.6666
$ws.Range("H64").Value = 5099.8
$ws.Range("I64").Value = 3499
$ws.Range("K64").Value = 3499
$ws.Range("M64").Value = -3251
$ws.Range("H67").Value = 5099.8
$ws.Range("I67").Value = 3499
$ws.Range("K67").Value = 3499
$ws.Range("M67").Value = -2641
$ws.Range("H88").Value = 2297.25
$ws.Range("I88").Value = 2199
$ws.Range("J88").Value = 2330
$ws.Range("K88").Value = 2199
$ws.Range("L88").Value = 2330
$ws.Range("M88").Value = -1793
$ws.Range("N88").Value = -3142
$ws.Range("H91").Value = 2297.25
$ws.Range("I91").Value = 2199
$ws.Range("J91").Value = 2330
$ws.Range("K91").Value = 2199
$ws.Range("L91").Value = 2330
$ws.Range("M91").Value = -795
$ws.Range("N91").Value = -5138
$ws.Range("H116").Value = 3582.8
$ws.Range("I116").Value = 3331.5
$ws.Range("J116").Value = 3750.3333
$ws.Range("K116").Value = 3331.5
$ws.Range("L116").Value = 3750.3333
$ws.Range("M116").Value = 110.5
$ws.Range("N116").Value = -10634.3333
$ws.Range("H137").Value = 5149.6665
$ws.Range("I137").Value = 1633.3334
$ws.Range("K137").Value = 4900.0002
$ws.Range("M137").Value = -2350.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3149.3333
$ws.Range("I2").Value = 2212.25
$ws.Range("J2").Value = 3899
$ws.Range("K2").Value = 2212.25
$ws.Range("L2").Value = 3899
$ws.Range("M2").Value = -2099.25
$ws.Range("N2").Value = -4125
$ws.Range("H88").Value = 2097
$ws.Range("I88").Value = 2206
$ws.Range("J88").Value = 2065.8572
$ws.Range("K88").Value = 2206
$ws.Range("L88").Value = 2065.8572
$ws.Range("M88").Value = -1800
$ws.Range("N88").Value = -2877.8572
$ws.Range("H91").Value = 2097
$ws.Range("I91").Value = 2206
$ws.Range("J91").Value = 2065.8572
$ws.Range("K91").Value = 2206
$ws.Range("L91").Value = 2065.8572
$ws.Range("M91").Value = -802
$ws.Range("N91").Value = -4873.8572
$ws.Range("H102").Value = 1389
$ws.Range("I102").Value = 1389
$ws.Range("K102").Value = 1389
$ws.Range("M102").Value = 233
$ws.Range("H110").Value = 1684.2222
$ws.Range("I110").Value = 1684.2222
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1684.2222
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 360.7778000000001
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 3149.3333
$ws.Range("I116").Value = 2212.25
$ws.Range("J116").Value = 3899
$ws.Range("K116").Value = 2212.25
$ws.Range("L116").Value = 3899
$ws.Range("M116").Value = 81.75
$ws.Range("N116").Value = -8487

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3149.3333
$ws.Range("I3").Value = 2212.25
$ws.Range("J3").Value = 3899
$ws.Range("K3").Value = 2212.25
$ws.Range("L3").Value = 3899
$ws.Range("M3").Value = -2098.25
$ws.Range("N3").Value = -4127
$ws.Range("H16").Value = 17330
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 17330
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 17330
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -17670
$ws.Range("H94").Value = 726.3333
$ws.Range("I94").Value = 829.6
$ws.Range("J94").Value = 210
$ws.Range("K94").Value = 829.6
$ws.Range("L94").Value = 210
$ws.Range("M94").Value = -378.6
$ws.Range("N94").Value = -1112
$ws.Range("H134").Value = 1868.4
$ws.Range("I134").Value = 1868.4
$ws.Range("K134").Value = 5605.200000000001
$ws.Range("M134").Value = -3070.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2659.1667
$ws.Range("I12").Value = 1227.5
$ws.Range("J12").Value = 3375
$ws.Range("K12").Value = 1227.5
$ws.Range("L12").Value = 3375
$ws.Range("M12").Value = -1057.5
$ws.Range("N12").Value = -3715
$ws.Range("H13").Value = 93342.5
$ws.Range("I13").Value = 185785
$ws.Range("J13").Value = 900
$ws.Range("K13").Value = 185785
$ws.Range("L13").Value = 900
$ws.Range("M13").Value = -185646
$ws.Range("N13").Value = -1178
$ws.Range("H134").Value = 1724
$ws.Range("I134").Value = 948
$ws.Range("K134").Value = 2844
$ws.Range("M134").Value = -309

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 200.3077
$ws.Range("I6").Value = 200.36363
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 601.0908899999999
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -488.0908899999999
$ws.Range("N6").Value = -826
$ws.Range("H81").Value = 2131.3333
$ws.Range("J81").Value = 2131.3333
$ws.Range("L81").Value = 6393.999899999999
$ws.Range("N81").Value = -8639.999899999999
$ws.Range("H82").Value = 1800
$ws.Range("I82").Value = 1800
$ws.Range("K82").Value = 5400
$ws.Range("M82").Value = -4994
$ws.Range("H84").Value = 2131.3333
$ws.Range("J84").Value = 2131.3333
$ws.Range("L84").Value = 19181.9997
$ws.Range("N84").Value = -30413.9997
$ws.Range("H85").Value = 1800
$ws.Range("I85").Value = 1800
$ws.Range("K85").Value = 5400
$ws.Range("M85").Value = -3996
$ws.Range("H109").Value = 1955.2
$ws.Range("I109").Value = 1744
$ws.Range("K109").Value = 5232
$ws.Range("M109").Value = -4192
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3826.8572
$ws.Range("I70").Value = 3947.25
$ws.Range("J70").Value = 3666.3333
$ws.Range("K70").Value = 3947.25
$ws.Range("L70").Value = 3666.3333
$ws.Range("M70").Value = -3677.25
$ws.Range("N70").Value = -4206.3333
$ws.Range("H73").Value = 3826.8572
$ws.Range("I73").Value = 3947.25
$ws.Range("J73").Value = 3666.3333
$ws.Range("K73").Value = 3947.25
$ws.Range("L73").Value = 3666.3333
$ws.Range("M73").Value = -3011.25
$ws.Range("N73").Value = -5538.3333
$ws.Range("H80").Value = 3228.4285
$ws.Range("I80").Value = 2924.75
$ws.Range("J80").Value = 3633.3333
$ws.Range("K80").Value = 2924.75
$ws.Range("L80").Value = 3633.3333
$ws.Range("M80").Value = -1926.75
$ws.Range("N80").Value = -5629.3333
$ws.Range("H83").Value = 3228.4285
$ws.Range("I83").Value = 2924.75
$ws.Range("J83").Value = 3633.3333
$ws.Range("K83").Value = 14623.75
$ws.Range("L83").Value = 18166.6665
$ws.Range("M83").Value = -9631.75
$ws.Range("N83").Value = -28150.6665
$ws.Range("H102").Value = 2288.5557
$ws.Range("I102").Value = 2296.875
$ws.Range("J102").Value = 2222
$ws.Range("K102").Value = 2296.875
$ws.Range("L102").Value = 2222
$ws.Range("M102").Value = -674.875
$ws.Range("N102").Value = -5466
$ws.Range("H132").Value = 5006.769
$ws.Range("I132").Value = 5006.769
$ws.Range("K132").Value = 15020.307
$ws.Range("M132").Value = -12490.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 999999.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 9199.200000000001
$ws.Range("J7").Value = 9199.200000000001
$ws.Range("L7").Value = 9199.200000000001
$ws.Range("N7").Value = -9423.200000000001
$ws.Range("H22").Value = 1959.2273
$ws.Range("J22").Value = 2729.3333
$ws.Range("L22").Value = 2729.3333
$ws.Range("N22").Value = -3319.3333
$ws.Range("H27").Value = 1959.2273
$ws.Range("J27").Value = 2729.3333
$ws.Range("L27").Value = 2729.3333
$ws.Range("N27").Value = -2943.3333
$ws.Range("H28").Value = 999999.5
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 999999.5
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H46").Value = 2423
$ws.Range("I46").Value = 2016
$ws.Range("K46").Value = 2016
$ws.Range("M46").Value = -1828
$ws.Range("H100").Value = 4472.4546
$ws.Range("I100").Value = 4472.4546
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4472.4546
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3931.4546
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 9199.200000000001
$ws.Range("J126").Value = 9199.200000000001
$ws.Range("L126").Value = 27597.6
$ws.Range("N126").Value = -32537.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 10000
$ws.Range("K24").Value = 10000
$ws.Range("M24").Value = -9770
$ws.Range("H136").Value = 776.36365
$ws.Range("I136").Value = 699.5
$ws.Range("J136").Value = 868.6
$ws.Range("K136").Value = 2098.5
$ws.Range("L136").Value = 2605.8
$ws.Range("M136").Value = 451.5
$ws.Range("N136").Value = -7705.8
